# "Generate Report for Handback"
#
# The two handoff files in this localization-status report have swapped
# places: the file that used to be "Ready for handoff" (61b896cf-...) has
# now been handed back (in sync with en-US), and its handback timestamps
# are filled in. The per-file rows on the language sheets ("zh-cn",
# "de-de") are re-ordered so the 61b896cf file is row 2 and the
# e0c3cbec file is row 3 (previously the opposite), and the "Overview"
# sheet is updated to match.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# Overview sheet
# ---------------------------------------------------------------------
$ov = $wb.Worksheets.Item("Overview")

$ov.Range("A2").Value2 = "61b896cf-cc6b-4613-bae6-25589e9c641c.md"
$ov.Range("B2").Value2 = "Handed back: in sync with en-US"
$ov.Range("C2").Value2 = "Handed back: in sync with en-US"
$ov.Range("D2").Value2 = "2016-03-23 22:50:31"

$ov.Range("A3").Value2 = "e0c3cbec-fd90-4bf9-b4dc-a0f9ed3e67c6.md"
$ov.Range("B3").Value2 = "Handed back: in sync with en-US"
$ov.Range("C3").Value2 = "Handed back: in sync with en-US"
$ov.Range("D3").Value2 = "2016-03-23 22:48:57"

# hyperlink display text follows the same swap (the underlying rIds /
# target URLs stay attached to the same cells)
$ov.Hyperlinks.Item(1).TextToDisplay = "61b896cf-cc6b-4613-bae6-25589e9c641c.md"
$ov.Hyperlinks.Item(2).TextToDisplay = "e0c3cbec-fd90-4bf9-b4dc-a0f9ed3e67c6.md"

# ---------------------------------------------------------------------
# zh-cn sheet
# ---------------------------------------------------------------------
$zh = $wb.Worksheets.Item("zh-cn")

$zh.Range("A2").Value2 = "61b896cf-cc6b-4613-bae6-25589e9c641c.md"
$zh.Range("B2").Value2 = ".md"
$zh.Range("C2").Value2 = "Handed back: in sync with en-US"
$zh.Range("D2").Value2 = "61b896cf-cc6b-4613-bae6-25589e9c641c.0cb423db10d2ca3cac4e4e2e5696829bdf7b154d.zh-cn.xlf"
$zh.Range("E2").Value2 = "2016-03-23 22:50:27"
$zh.Range("F2").Value2 = "61b896cf-cc6b-4613-bae6-25589e9c641c.md"
$zh.Range("G2").Value2 = "61b896cf-cc6b-4613-bae6-25589e9c641c.0cb423db10d2ca3cac4e4e2e5696829bdf7b154d.zh-cn.xlf"
$zh.Range("H2").Value2 = "2016-03-23 22:50:51"
$zh.Range("J2").Value2 = "Include"

$zh.Range("A3").Value2 = "e0c3cbec-fd90-4bf9-b4dc-a0f9ed3e67c6.md"
$zh.Range("B3").Value2 = ".md"
$zh.Range("C3").Value2 = "Handed back: in sync with en-US"
$zh.Range("D3").Value2 = "e0c3cbec-fd90-4bf9-b4dc-a0f9ed3e67c6.f46fd9bbdb5bce68e26b2f9491a78b463d29c64c.zh-cn.xlf"
$zh.Range("E3").Value2 = "2016-03-23 22:48:53"
$zh.Range("F3").Value2 = "e0c3cbec-fd90-4bf9-b4dc-a0f9ed3e67c6.md"
$zh.Range("G3").Value2 = "e0c3cbec-fd90-4bf9-b4dc-a0f9ed3e67c6.f46fd9bbdb5bce68e26b2f9491a78b463d29c64c.zh-cn.xlf"
$zh.Range("H3").Value2 = "2016-03-23 22:49:31"
$zh.Range("J3").Value2 = "Include"

$zh.Hyperlinks.Item(1).TextToDisplay = "61b896cf-cc6b-4613-bae6-25589e9c641c.md"
$zh.Hyperlinks.Item(2).TextToDisplay = "61b896cf-cc6b-4613-bae6-25589e9c641c.0cb423db10d2ca3cac4e4e2e5696829bdf7b154d.zh-cn.xlf"
$zh.Hyperlinks.Item(3).TextToDisplay = "61b896cf-cc6b-4613-bae6-25589e9c641c.md"
$zh.Hyperlinks.Item(4).TextToDisplay = "61b896cf-cc6b-4613-bae6-25589e9c641c.0cb423db10d2ca3cac4e4e2e5696829bdf7b154d.zh-cn.xlf"
$zh.Hyperlinks.Item(5).TextToDisplay = "e0c3cbec-fd90-4bf9-b4dc-a0f9ed3e67c6.md"
$zh.Hyperlinks.Item(6).TextToDisplay = "e0c3cbec-fd90-4bf9-b4dc-a0f9ed3e67c6.f46fd9bbdb5bce68e26b2f9491a78b463d29c64c.zh-cn.xlf"
$zh.Hyperlinks.Item(7).TextToDisplay = "e0c3cbec-fd90-4bf9-b4dc-a0f9ed3e67c6.md"
$zh.Hyperlinks.Item(8).TextToDisplay = "e0c3cbec-fd90-4bf9-b4dc-a0f9ed3e67c6.f46fd9bbdb5bce68e26b2f9491a78b463d29c64c.zh-cn.xlf"

# ---------------------------------------------------------------------
# de-de sheet
# ---------------------------------------------------------------------
$de = $wb.Worksheets.Item("de-de")

$de.Range("A2").Value2 = "61b896cf-cc6b-4613-bae6-25589e9c641c.md"
$de.Range("B2").Value2 = ".md"
$de.Range("C2").Value2 = "Handed back: in sync with en-US"
$de.Range("D2").Value2 = "61b896cf-cc6b-4613-bae6-25589e9c641c.0cb423db10d2ca3cac4e4e2e5696829bdf7b154d.de-de.xlf"
$de.Range("E2").Value2 = "2016-03-23 22:50:31"
$de.Range("F2").Value2 = "61b896cf-cc6b-4613-bae6-25589e9c641c.md"
$de.Range("G2").Value2 = "61b896cf-cc6b-4613-bae6-25589e9c641c.0cb423db10d2ca3cac4e4e2e5696829bdf7b154d.de-de.xlf"
$de.Range("H2").Value2 = "2016-03-23 22:50:57"
$de.Range("J2").Value2 = "Include"

$de.Range("A3").Value2 = "e0c3cbec-fd90-4bf9-b4dc-a0f9ed3e67c6.md"
$de.Range("B3").Value2 = ".md"
$de.Range("C3").Value2 = "Handed back: in sync with en-US"
$de.Range("D3").Value2 = "e0c3cbec-fd90-4bf9-b4dc-a0f9ed3e67c6.f46fd9bbdb5bce68e26b2f9491a78b463d29c64c.de-de.xlf"
$de.Range("E3").Value2 = "2016-03-23 22:48:57"
$de.Range("F3").Value2 = "e0c3cbec-fd90-4bf9-b4dc-a0f9ed3e67c6.md"
$de.Range("G3").Value2 = "e0c3cbec-fd90-4bf9-b4dc-a0f9ed3e67c6.f46fd9bbdb5bce68e26b2f9491a78b463d29c64c.de-de.xlf"
$de.Range("H3").Value2 = "2016-03-23 22:49:40"
$de.Range("J3").Value2 = "Include"

$de.Hyperlinks.Item(1).TextToDisplay = "61b896cf-cc6b-4613-bae6-25589e9c641c.md"
$de.Hyperlinks.Item(2).TextToDisplay = "61b896cf-cc6b-4613-bae6-25589e9c641c.0cb423db10d2ca3cac4e4e2e5696829bdf7b154d.de-de.xlf"
$de.Hyperlinks.Item(3).TextToDisplay = "61b896cf-cc6b-4613-bae6-25589e9c641c.md"
$de.Hyperlinks.Item(4).TextToDisplay = "61b896cf-cc6b-4613-bae6-25589e9c641c.0cb423db10d2ca3cac4e4e2e5696829bdf7b154d.de-de.xlf"
$de.Hyperlinks.Item(5).TextToDisplay = "e0c3cbec-fd90-4bf9-b4dc-a0f9ed3e67c6.md"
$de.Hyperlinks.Item(6).TextToDisplay = "e0c3cbec-fd90-4bf9-b4dc-a0f9ed3e67c6.f46fd9bbdb5bce68e26b2f9491a78b463d29c64c.de-de.xlf"
$de.Hyperlinks.Item(7).TextToDisplay = "e0c3cbec-fd90-4bf9-b4dc-a0f9ed3e67c6.md"
$de.Hyperlinks.Item(8).TextToDisplay = "e0c3cbec-fd90-4bf9-b4dc-a0f9ed3e67c6.f46fd9bbdb5bce68e26b2f9491a78b463d29c64c.de-de.xlf"
